# Updates per-row leve-profit figures (columns H-N) across several worksheets,
# reflecting refreshed market-board price data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1566.375
$ws.Range("I19").Value = 2299.8
$ws.Range("K19").Value = 2299.8
$ws.Range("M19").Value = -2124.8

# Row 70
$ws.Range("H70").Value = 2750
$ws.Range("I70").Value = 2500
$ws.Range("J70").Value = 2833.3333
$ws.Range("K70").Value = 7500
$ws.Range("L70").Value = 8499.999899999999
$ws.Range("M70").Value = -7230
$ws.Range("N70").Value = -9039.999899999999

# Row 73
$ws.Range("H73").Value = 2750
$ws.Range("I73").Value = 2500
$ws.Range("J73").Value = 2833.3333
$ws.Range("K73").Value = 7500
$ws.Range("L73").Value = 8499.999899999999
$ws.Range("M73").Value = -6564
$ws.Range("N73").Value = -10371.9999

# Row 98
$ws.Range("H98").Value = 979.05884
$ws.Range("I98").Value = 1129.4
$ws.Range("J98").Value = 764.2857
$ws.Range("K98").Value = 1129.4
$ws.Range("L98").Value = 764.2857
$ws.Range("M98").Value = 368.5999999999999
$ws.Range("N98").Value = -3760.2857

# Row 122
$ws.Range("H122").Value = 979.05884
$ws.Range("I122").Value = 1129.4
$ws.Range("J122").Value = 764.2857
$ws.Range("K122").Value = 3388.2
$ws.Range("L122").Value = 2292.8571
$ws.Range("M122").Value = -938.2000000000003
$ws.Range("N122").Value = -7192.8571

# Row 137
$ws.Range("H137").Value = 2209.5217
$ws.Range("I137").Value = 1632.579
$ws.Range("K137").Value = 4897.737
$ws.Range("M137").Value = -2347.737

# Row 138
$ws.Range("H138").Value = 3731.9285
$ws.Range("J138").Value = 5149.75
$ws.Range("L138").Value = 15449.25
$ws.Range("N138").Value = -25729.25

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3910.7646
$ws.Range("I32").Value = 2676.3547
$ws.Range("J32").Value = 16666.334
$ws.Range("K32").Value = 2676.3547
$ws.Range("L32").Value = 16666.334
$ws.Range("M32").Value = -2389.3547
$ws.Range("N32").Value = -17240.334

# Row 102
$ws.Range("H102").Value = 2060.375
$ws.Range("I102").Value = 1247.5
$ws.Range("K102").Value = 1247.5
$ws.Range("M102").Value = 374.5

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 411
$ws.Range("I22").Value = 411
$ws.Range("K22").Value = 411
$ws.Range("M22").Value = -238

# Row 80
$ws.Range("H80").Value = 402.23077
$ws.Range("J80").Value = 336
$ws.Range("L80").Value = 336
$ws.Range("N80").Value = -2332

# Row 83
$ws.Range("H83").Value = 402.23077
$ws.Range("J83").Value = 336
$ws.Range("L83").Value = 1680
$ws.Range("N83").Value = -11664

# Row 94
$ws.Range("H94").Value = 4750.778
$ws.Range("I94").Value = 4551.4
$ws.Range("K94").Value = 4551.4
$ws.Range("M94").Value = -4100.4

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2114
$ws.Range("J22").Value = 3953.3333
$ws.Range("L22").Value = 3953.3333
$ws.Range("N22").Value = -4653.3333

# Row 62
$ws.Range("H62").Value = 2666.3333
$ws.Range("I62").Value = 2499.5
$ws.Range("K62").Value = 2499.5
$ws.Range("M62").Value = -1875.5

# Row 65
$ws.Range("H65").Value = 2666.3333
$ws.Range("I65").Value = 2499.5
$ws.Range("K65").Value = 12497.5
$ws.Range("M65").Value = -9377.5

# Row 99
$ws.Range("H99").Value = 4620.15
$ws.Range("I99").Value = 5208.9165
$ws.Range("J99").Value = 3737
$ws.Range("K99").Value = 5208.9165
$ws.Range("L99").Value = 3737
$ws.Range("M99").Value = -3710.9165
$ws.Range("N99").Value = -6733

# Row 122
$ws.Range("H122").Value = 1797.7646
$ws.Range("I122").Value = 1006.375
$ws.Range("K122").Value = 3019.125
$ws.Range("M122").Value = -569.125

# Row 126
$ws.Range("H126").Value = 4620.15
$ws.Range("I126").Value = 5208.9165
$ws.Range("J126").Value = 3737
$ws.Range("K126").Value = 15626.7495
$ws.Range("L126").Value = 11211
$ws.Range("M126").Value = -13156.7495
$ws.Range("N126").Value = -16151

# Row 132
$ws.Range("H132").Value = 1808.5
$ws.Range("I132").Value = 1919.8
$ws.Range("K132").Value = 5759.4
$ws.Range("M132").Value = -3229.4

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 732
$ws.Range("I55").Value = 98.25
$ws.Range("J55").Value = 1999.5
$ws.Range("K55").Value = 294.75
$ws.Range("L55").Value = 5998.5
$ws.Range("M55").Value = -117.75
$ws.Range("N55").Value = -6352.5

$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 95
$ws.Range("H95").Value = 24614.666
$ws.Range("J95").Value = 25172
$ws.Range("L95").Value = 25172
$ws.Range("N95").Value = -30664

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2397.5
$ws.Range("I7").Value = 1998.375
$ws.Range("K7").Value = 1998.375
$ws.Range("M7").Value = -1886.375

# Row 16
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1330

# Row 40
$ws.Range("H40").Value = 4412.4585
$ws.Range("I40").Value = 4245.476
$ws.Range("K40").Value = 4245.476
$ws.Range("M40").Value = -4109.476

# Row 46
$ws.Range("H46").Value = 1534.0968
$ws.Range("J46").Value = 2064.182
$ws.Range("L46").Value = 2064.182
$ws.Range("N46").Value = -2440.182

# Row 68
$ws.Range("H68").Value = 2049.5
$ws.Range("I68").Value = 1398.3334
$ws.Range("K68").Value = 1398.3334
$ws.Range("M68").Value = -649.3334

# Row 71
$ws.Range("H71").Value = 2049.5
$ws.Range("I71").Value = 1398.3334
$ws.Range("K71").Value = 6991.666999999999
$ws.Range("M71").Value = -3247.666999999999

# Row 82
$ws.Range("H82").Value = 992.3333
$ws.Range("I82").Value = 992.3333
$ws.Range("K82").Value = 992.3333
$ws.Range("M82").Value = -631.3333

# Row 85
$ws.Range("H85").Value = 992.3333
$ws.Range("I85").Value = 992.3333
$ws.Range("K85").Value = 992.3333
$ws.Range("M85").Value = 255.6667

# Row 126
$ws.Range("H126").Value = 2397.5
$ws.Range("I126").Value = 1998.375
$ws.Range("K126").Value = 5995.125
$ws.Range("M126").Value = -3525.125

# Row 130
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 100000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 100000
$ws.Range("N130").Value = -110040
$ws.Range("M130").ClearContents()

# Row 132
$ws.Range("H132").Value = 5180.684
$ws.Range("J132").Value = 3833
$ws.Range("L132").Value = 11499
$ws.Range("N132").Value = -16559

# Row 136
$ws.Range("H136").Value = 4256.7856
$ws.Range("J136").Value = 4966.6665
$ws.Range("L136").Value = 14899.9995
$ws.Range("N136").Value = -19999.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5397
$ws.Range("I62").Value = 5330.3335
$ws.Range("J62").Value = 5497
$ws.Range("K62").Value = 5330.3335
$ws.Range("L62").Value = 5497
$ws.Range("M62").Value = -4706.3335
$ws.Range("N62").Value = -6745

# Row 65
$ws.Range("H65").Value = 5397
$ws.Range("I65").Value = 5330.3335
$ws.Range("J65").Value = 5497
$ws.Range("K65").Value = 26651.6675
$ws.Range("L65").Value = 27485
$ws.Range("M65").Value = -23531.6675
$ws.Range("N65").Value = -33725

# Row 81
$ws.Range("H81").Value = 7124.5
$ws.Range("I81").Value = 5000
$ws.Range("J81").Value = 7832.6665
$ws.Range("K81").Value = 10000
$ws.Range("L81").Value = 15665.333
$ws.Range("M81").Value = -8939
$ws.Range("N81").Value = -17787.333

# Row 84
$ws.Range("H84").Value = 7124.5
$ws.Range("I84").Value = 5000
$ws.Range("J84").Value = 7832.6665
$ws.Range("K84").Value = 50000
$ws.Range("L84").Value = 78326.66500000001
$ws.Range("M84").Value = -44696
$ws.Range("N84").Value = -88934.66500000001

# Row 100
$ws.Range("H100").Value = 4981233.5
$ws.Range("I100").Value = 9960076
$ws.Range("J100").Value = 2391.4285
$ws.Range("K100").Value = 19920152
$ws.Range("L100").Value = 4782.857
$ws.Range("M100").Value = -19919611
$ws.Range("N100").Value = -5864.857

# Row 126
$ws.Range("H126").Value = 2808.8333
$ws.Range("I126").Value = 2941.4546
$ws.Range("K126").Value = 8824.363799999999
$ws.Range("M126").Value = -6354.363799999999

# Row 132
$ws.Range("H132").Value = 765.1429000000001
$ws.Range("I132").Value = 553.4
$ws.Range("K132").Value = 1660.2
$ws.Range("M132").Value = 869.8000000000002
